$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Supplier row (row 9): "Confirm Filtered Record (Without Duplication)" (E)
# and "Progress" (F) both change from 435 to 1366.
$ws.Range("E9").Value = 1366
$ws.Range("F9").Value = 1366

# Re-enter the ratio formulas across their full blocks so Excel collapses
# them into shared formulas (matches how the workbook was actually saved).
$ws.Range("G5:G12").Formula = "=IF(D5>0, F5/D5, """")"
$ws.Range("H5:H14").Formula = "=IF(E5>0, F5/E5, """")"

# Reflect the last selected/active cell as seen in the saved file.
$ws.Range("F9").Select()
